$wb = $excel.ActiveWorkbook

# --- Update findCarModelAndPriceTest sheet (swap Hynundai/Toyota rows, selection) ---
$carSheet = $wb.Worksheets.Item("findCarModelAndPriceTest")
$carSheet.Range("A4").Value = "Toyota"
$carSheet.Range("A5").Value = "Hynundai"
$carSheet.Range("B1:B5").Select() | Out-Null

# --- Add the new homePageTest sheet at the end ---
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "homePageTest"
$newSheet.Range("A1").Value = "browserName"
$newSheet.Range("A2").Value = "chrome"
$newSheet.Range("C6").Select() | Out-Null
